# Update fruit_name values on Sheet1 (B3: Apple -> Iphone, B5: Banana -> Republic)
$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Sheet1")
$ws1.Range("B3").Value = "Iphone"
$ws1.Range("B5").Value = "Republic"

# Restore the active selection on Sheet1 to match the saved view state
$ws1.Range("G14").Select() | Out-Null
